# Apply updated profit/price figures across all 8 job sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 915
$ws.Range("I6").Value = 451.66666
$ws.Range("K6").Value = 1354.99998
$ws.Range("M6").Value = -1242.99998

$ws.Range("H41").Value = 527.3333
$ws.Range("I41").Value = 425.42856
$ws.Range("J41").Value = 884
$ws.Range("K41").Value = 425.42856
$ws.Range("L41").Value = 884
$ws.Range("M41").Value = 14.57144
$ws.Range("N41").Value = -1764

$ws.Range("H70").Value = 3129.5557
$ws.Range("I70").Value = 2790
$ws.Range("J70").Value = 3139.257
$ws.Range("K70").Value = 8370
$ws.Range("L70").Value = 9417.771000000001
$ws.Range("M70").Value = -8100
$ws.Range("N70").Value = -9957.771000000001

$ws.Range("H73").Value = 3129.5557
$ws.Range("I73").Value = 2790
$ws.Range("J73").Value = 3139.257
$ws.Range("K73").Value = 8370
$ws.Range("L73").Value = 9417.771000000001
$ws.Range("M73").Value = -7434
$ws.Range("N73").Value = -11289.771

$ws.Range("H138").Value = 2960.1072
$ws.Range("J138").Value = 3489.1428
$ws.Range("L138").Value = 10467.4284
$ws.Range("N138").Value = -20747.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 16502
$ws.Range("I13").Value = 20000
$ws.Range("J13").Value = 13004
$ws.Range("K13").Value = 20000
$ws.Range("L13").Value = 13004
$ws.Range("M13").Value = -19856
$ws.Range("N13").Value = -13292

$ws.Range("H31").Value = 6020.778
$ws.Range("I31").Value = 6020.778
$ws.Range("K31").Value = 6020.778
$ws.Range("M31").Value = -5726.778

$ws.Range("H32").Value = 17177.4
$ws.Range("I32").Value = 19482.26
$ws.Range("K32").Value = 19482.26
$ws.Range("M32").Value = -19195.26

$ws.Range("H45").Value = 3884
$ws.Range("I45").Value = 2208.5
$ws.Range("J45").Value = 4721.75
$ws.Range("K45").Value = 2208.5
$ws.Range("L45").Value = 4721.75
$ws.Range("M45").Value = -1831.5
$ws.Range("N45").Value = -5475.75

$ws.Range("H74").Value = 1728.4333
$ws.Range("I74").Value = 1324.5769
$ws.Range("J74").Value = 4353.5
$ws.Range("K74").Value = 1324.5769
$ws.Range("L74").Value = 4353.5
$ws.Range("M74").Value = -450.5769
$ws.Range("N74").Value = -6101.5

$ws.Range("H77").Value = 1728.4333
$ws.Range("I77").Value = 1324.5769
$ws.Range("J77").Value = 4353.5
$ws.Range("K77").Value = 6622.8845
$ws.Range("L77").Value = 21767.5
$ws.Range("M77").Value = -2254.8845
$ws.Range("N77").Value = -30503.5

$ws.Range("H119").Value = 32832.832
$ws.Range("I119").Value = 19999
$ws.Range("J119").Value = 45666.668
$ws.Range("K119").Value = 19999
$ws.Range("L119").Value = 45666.668
$ws.Range("M119").Value = -15161
$ws.Range("N119").Value = -55342.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 67775.375
$ws.Range("I99").Value = 103459
$ws.Range("J99").Value = 8302.666999999999
$ws.Range("K99").Value = 103459
$ws.Range("L99").Value = 8302.666999999999
$ws.Range("M99").Value = -101961
$ws.Range("N99").Value = -11298.667

$ws.Range("H107").Value = 3519.5
$ws.Range("I107").Value = 3197.5
$ws.Range("J107").Value = 3841.5
$ws.Range("K107").Value = 3197.5
$ws.Range("L107").Value = 3841.5
$ws.Range("M107").Value = -1277.5
$ws.Range("N107").Value = -7681.5

$ws.Range("H123").Value = 57389.5
$ws.Range("J123").Value = 57389.5
$ws.Range("L123").Value = 57389.5
$ws.Range("N123").Value = -67189.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3428.353
$ws.Range("I16").Value = 3125.5
$ws.Range("J16").Value = 3861
$ws.Range("K16").Value = 3125.5
$ws.Range("L16").Value = 3861
$ws.Range("M16").Value = -2838.5
$ws.Range("N16").Value = -4435

$ws.Range("H31").Value = 1347.3704
$ws.Range("I31").Value = 1378.9584
$ws.Range("K31").Value = 1378.9584
$ws.Range("M31").Value = -1083.9584

$ws.Range("H34").Value = 1347.3704
$ws.Range("I34").Value = 1378.9584
$ws.Range("K34").Value = 1378.9584
$ws.Range("M34").Value = -1176.9584

$ws.Range("H113").Value = 3428.353
$ws.Range("I113").Value = 3125.5
$ws.Range("J113").Value = 3861
$ws.Range("K113").Value = 3125.5
$ws.Range("L113").Value = 3861
$ws.Range("M113").Value = -955.5
$ws.Range("N113").Value = -8201

$ws.Range("H132").Value = 1928.8334
$ws.Range("I132").Value = 1928.8334
$ws.Range("K132").Value = 5786.5002
$ws.Range("M132").Value = -3256.5002

$ws.Range("H134").Value = 127322.5
$ws.Range("I134").Value = 127322.5
$ws.Range("K134").Value = 381967.5
$ws.Range("M134").Value = -379432.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 467
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

$ws.Range("H127").Value = 7630.5
$ws.Range("J127").Value = 7630.5
$ws.Range("L127").Value = 22891.5
$ws.Range("N127").Value = -32811.5

$ws.Range("H137").Value = 2810.7646
$ws.Range("I137").Value = 1999.4546
$ws.Range("K137").Value = 5998.3638
$ws.Range("M137").Value = -898.3638000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 29449.857
$ws.Range("I14").Value = 16229.8
$ws.Range("J14").Value = 62500
$ws.Range("K14").Value = 16229.8
$ws.Range("L14").Value = 62500
$ws.Range("M14").Value = -16061.8
$ws.Range("N14").Value = -62836

$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()

$ws.Range("H126").Value = 5208.6665
$ws.Range("J126").Value = 4724.75
$ws.Range("L126").Value = 14174.25
$ws.Range("N126").Value = -19114.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 49999.668
$ws.Range("J14").Value = 49999.668
$ws.Range("L14").Value = 49999.668
$ws.Range("N14").Value = -50343.668

$ws.Range("H16").Value = 2302.9167
$ws.Range("I16").Value = 2136.5264
$ws.Range("J16").Value = 2935.2
$ws.Range("K16").Value = 2136.5264
$ws.Range("L16").Value = 2935.2
$ws.Range("M16").Value = -1966.5264
$ws.Range("N16").Value = -3275.2

$ws.Range("H55").Value = 483
$ws.Range("I55").Value = 385
$ws.Range("J55").Value = 875
$ws.Range("K55").Value = 385
$ws.Range("L55").Value = 875
$ws.Range("M55").Value = -212
$ws.Range("N55").Value = -1221

$ws.Range("H61").Value = 7400.212
$ws.Range("I61").Value = 6569.6
$ws.Range("J61").Value = 9995.875
$ws.Range("K61").Value = 6569.6
$ws.Range("L61").Value = 9995.875
$ws.Range("M61").Value = -6367.6
$ws.Range("N61").Value = -10399.875

$ws.Range("H93").Value = 1995.5416
$ws.Range("I93").Value = 1840.0625
$ws.Range("J93").Value = 2306.5
$ws.Range("K93").Value = 1840.0625
$ws.Range("L93").Value = 2306.5
$ws.Range("M93").Value = -592.0625
$ws.Range("N93").Value = -4802.5

$ws.Range("H113").Value = 7400.212
$ws.Range("I113").Value = 6569.6
$ws.Range("J113").Value = 9995.875
$ws.Range("K113").Value = 6569.6
$ws.Range("L113").Value = 9995.875
$ws.Range("M113").Value = -4399.6
$ws.Range("N113").Value = -14335.875

$ws.Range("H119").Value = 19999.334
$ws.Range("I119").Value = 19999
$ws.Range("J119").Value = 20000
$ws.Range("K119").Value = 19999
$ws.Range("L119").Value = 20000
$ws.Range("M119").Value = -15161
$ws.Range("N119").Value = -29676

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 35996.777
$ws.Range("J45").Value = 35996.777
$ws.Range("L45").Value = 35996.777
$ws.Range("N45").Value = -36978.777

$ws.Range("H62").Value = 130750.375
$ws.Range("I62").Value = 3000
$ws.Range("K62").Value = 3000
$ws.Range("M62").Value = -2376

$ws.Range("H65").Value = 130750.375
$ws.Range("I65").Value = 3000
$ws.Range("K65").Value = 15000
$ws.Range("M65").Value = -11880

$ws.Range("H113").Value = 4592.75
$ws.Range("I113").Value = 3618.889
$ws.Range("K113").Value = 10856.667
$ws.Range("M113").Value = -8686.667000000001

$ws.Range("H119").Value = 60000
$ws.Range("J119").Value = 60000
$ws.Range("L119").Value = 60000
$ws.Range("N119").Value = -69676

$ws.Range("H122").Value = 5588
$ws.Range("J122").Value = 9964
$ws.Range("L122").Value = 29892
$ws.Range("N122").Value = -34792

$ws.Range("H132").Value = 34480.258
$ws.Range("I132").Value = 35396.266
$ws.Range("K132").Value = 106188.798
$ws.Range("M132").Value = -103658.798

$ws.Range("H136").Value = 3910.111
$ws.Range("I136").Value = 2576.4614
$ws.Range("K136").Value = 7729.3842
$ws.Range("M136").Value = -5179.3842

